$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date serial changes from 45406 (2024-04-24) to 45436 (2024-05-24)
$ws.Range("A1").Value = 45436

# D29: price changes from 960 to 2100
$ws.Range("D29").Value = 2100
